$wb = $excel.ActiveWorkbook

# --- Import existing generation capacity into "Info geografiche" (NORD row) ---
$wsGeo = $wb.Worksheets.Item("Info geografiche")
$wsGeo.Range("B2").Value = 9999
$wsGeo.Range("C2").Value = 9999
$wsGeo.Range("D2").Value = 0

# --- Fill the previously-blank "SUD " duplicate row on the capacities sheet ---
$wsCap = $wb.Worksheets.Item("Capacità di trasmissione MW")
$wsCap.Range("B9").Value = 0
$wsCap.Range("C9").Value = 0
$wsCap.Range("D9").Value = 0
$wsCap.Range("E9").Value = 0
$wsCap.Range("G9").Value = 0
$wsCap.Range("H9").Value = 0

# --- Re-assert the header / row-label formatting (bold, thin box border, centered) ---
$capHeader = $wsCap.Range("B1:H1")
$capHeader.Font.Bold = $true
$capHeader.Font.Name = "Aptos Narrow"
$capHeader.Font.Size = 11
$capHeader.Borders.LineStyle = 1
$capHeader.HorizontalAlignment = -4108
$capHeader.VerticalAlignment = -4160

$capRowLabels = $wsCap.Range("A2:A9")
$capRowLabels.Font.Bold = $true
$capRowLabels.Font.Name = "Aptos Narrow"
$capRowLabels.Font.Size = 11
$capRowLabels.Borders.LineStyle = 1
$capRowLabels.HorizontalAlignment = -4108
$capRowLabels.VerticalAlignment = -4160

# --- Switch the active tab back to "Info geografiche" with E5 selected ---
$wsGeo.Activate()
$wsGeo.Range("E5").Select()
